$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 149.4883017415314
$ws.Range("B3").Value = 144.7113422432715
$ws.Range("B4").Value = 147.27590589659
$ws.Range("B5").Value = 145.7066161003521
$ws.Range("B6").Value = 148.8324072647092
$ws.Range("B7").Value = 155.1171063241644
$ws.Range("B8").Value = 153.6704701299099
$ws.Range("B9").Value = 198.6334567207135
$ws.Range("B10").Value = 362.1433803996547
$ws.Range("B11").Value = 68.70107168834457
$ws.Range("B12").Value = 73.56977806884662
$ws.Range("B13").Value = 189.9157425227596
$ws.Range("B14").Value = 47.02370745953354
$ws.Range("B15").Value = 137.4376429993846
$ws.Range("B16").Value = 316.5360010704878
$ws.Range("B17").Value = 286.858166607642
$ws.Range("B18").Value = 149.182627445449
$ws.Range("B19").Value = 264.2715381946999
$ws.Range("B20").Value = 365.6947572331626
$ws.Range("B21").Value = 179.6626617176126
$ws.Range("B22").Value = 298.6928625621471
$ws.Range("B23").Value = 298.8082638992933
$ws.Range("B24").Value = 59.90295997935799
$ws.Range("B25").Value = 172.6149770563653
$ws.Range("B26").Value = 284.0734371414941
$ws.Range("B27").Value = 336.2616647642
$ws.Range("B28").Value = 311.2186250419591
$ws.Range("B29").Value = 310.0406837930459
$ws.Range("B30").Value = 138.4669318193171
$ws.Range("B31").Value = 154.5108555958332
$ws.Range("B32").Value = 139.5975827705442
$ws.Range("B33").Value = 74.28471843067288
$ws.Range("B34").Value = 147.2890190827495
